# Update the two TIAC document labels.
# Note: the document uses a non-breaking space (U+00A0) right before each
# colon (French typographic convention), so we build the search/replace
# strings explicitly with [char]0x00A0 rather than a plain space.

$d = $word.ActiveDocument
$nbsp = [char]0x00A0

# 1) "Date de réception à la DD(ETS)PP : ..." -> "Date de réception : ..."
$old1 = "Date de réception à la DD(ETS)PP" + $nbsp + ":"
$new1 = "Date de réception" + $nbsp + ":"
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false,
                         $true, 1, $false, $new1, 2)

# 2) "Suite donnée par la DD : ..." -> "Suite donnée : ..."
$old2 = "Suite donnée par la DD" + $nbsp + ":"
$new2 = "Suite donnée" + $nbsp + ":"
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false,
                         $true, 1, $false, $new2, 2)
